$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '28.636.19'
$ws.Range("D3").Value = '1.793.48'
$ws.Range("E3").Value = '  -2.03%  '
$ws.Range("E4").Value = '  +0.19%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '231.38'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -1.66%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.5887'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -2.45%  '
$ws.Range("E7").Value = '  +0.08%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.2764'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -1.01%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06742'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -4.43%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '23.13'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -1.79%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07522'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -1.82%  '
$ws.Range("D12").Value = '1.795.78'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '4.791'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -0.03%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.6137'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -2.42%  '
$ws.Range("D15").Value = '2.037.28'
$ws.Range("E15").Value = '  -2.06%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '75.24'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -4.88%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.000009016'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -8.76%  '
$ws.Range("D18").Value = '28.620.39'
$ws.Range("E18").Value = '  -2.31%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '5.456'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -6.74%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '1.005'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +0.10%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '209.83'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -6.53%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '11.46'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -2.08%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '6.810'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -2.81%  '
$ws.Range("E24").Value = '  +0.08%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '153.15'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -1.86%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '8.052'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +0.78%  '
$ws.Range("E27").Value = '  -3.61%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '16.40'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -1.17%  '
$ws.Range("E29").Value = '  -4.44%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.06088'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -5.09%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.421'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -1.92%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.810'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +0.19%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.784'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -1.47%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.737'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +0.25%  '
$ws.Range("E35").Value = '  -5.87%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.6403'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -1.02%  '
$ws.Range("E37").Value = '  -1.80%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.715'
$ws.Range("D38").Style = "Normal"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '6.409'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -2.17%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.01695'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -2.93%  '
$ws.Range("D41").Value = '1.142.07'
$ws.Range("E41").Value = '  -6.22%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.8801'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -2.26%  '
$ws.Range("E43").Value = '  +0.31%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '100.05'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -0.19%  '
$ws.Range("D45").Value = '1.945.70'
$ws.Range("E45").Value = '  -2.63%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '59.87'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -4.61%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.00000000112'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -3.87%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.583'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +0.01%  '
$ws.Range("B49").Value = 'Cronos'
$ws.Range("C49").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.05484'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -0.37%  '
$ws.Range("B50").Value = 'EnergySwap'
$ws.Range("C50").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '8.338'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -2.77%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.4483'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -1.53%  '
